# Updated symbol list on Sat Dec 17 11:34:52 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores its values as text in the source workbook
# (e.g. "238.37"), so every write into that column first forces the cell's
# number format to Text ("@") -- exactly like a human typing numbers into a
# pre-formatted text column in Excel -- to avoid Excel's automatic
# "look-like-a-number => convert to number" coercion, which would otherwise
# mangle values (lose trailing zeros, introduce float noise, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- simple price refreshes (rows 2-17) ---
Set-TextValue "D2"  "238.43"
Set-TextValue "D3"  "21.63"
Set-TextValue "D4"  "5.462"
Set-TextValue "D5"  "0.05644"
Set-TextValue "D6"  "6.488"
Set-TextValue "D9"  "0.7935"
Set-TextValue "D10" "0.1393"
Set-TextValue "D11" "0.07336"
Set-TextValue "D12" "0.03213"
Set-TextValue "D13" "0.02970"
Set-TextValue "D14" "0.09243"
Set-TextValue "D15" "0.001663"
Set-TextValue "D16" "3.257"
Set-TextValue "D17" "0.04771"

# --- rows 18-24: coin list re-ranked (rows shift up by one, "One" drops to
#     the bottom of this block with fresh data) ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.006207"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.005108"
$ws.Range("E19").Value = "18HotbitTokenHTB"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.001052"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.885"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D23" "2.198"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D24" "0.01172"
$ws.Range("E24").Value = "23OneONEBestin24h"

# --- more simple price refreshes (rows 40-41) ---
Set-TextValue "D40" "0.04133"
Set-TextValue "D41" "0.006955"

# --- rows 42-43: BKEXToken / CEJI swap places ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003501"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1040"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- row 44: price refresh + drop the "Bestin24h" suffix ---
Set-TextValue "D44" "0.009903"
$ws.Range("E44").Value = "43LocalTradersLCT"

# --- remaining simple price refreshes ---
Set-TextValue "D45" "0.00005444"
Set-TextValue "D47" "0.6755"
Set-TextValue "D48" "0.03730"
